# Update the intro slide's "HILT 2018" to "HILT 2019".
#
# Slide 1 has a body placeholder whose text is built from three runs:
#   "HILT " + "20" + "18"
# Only the final run's text ("18") needs to change to "19"; its run
# properties (size/color/lang) must be left untouched, so we target just
# the last two characters of the text range instead of overwriting the
# whole string.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item(2)          # "HILT 2018" placeholder (Shape 34)
$tr  = $shp.TextFrame.TextRange

# Characters(Start, Length) is 1-indexed; "HILT 2018" -> chars 8-9 are "18".
$yearDigits = $tr.Characters($tr.Text.Length - 1, 2)
$yearDigits.Text = "19"
